# wzrv/expdata/2012.xlsx -- "Add files via upload" re-upload edit
#
# The commit re-uploads the workbook with a few small, user-visible tweaks:
#   1. The shared string "CMS(e) " (trailing space) becomes "CMS(e)(2011)"
#      for every row of column L (cms label column).
#   2. The saved cursor/selection on Sheet1 moves from N22 to N7.
#   3. The saved window size grows (windowHeight 12300 -> 13590) and the
#      recorded absolute folder path in the x15ac:absPath hint changes --
#      both are host/window chrome metadata written by the real Excel
#      client on save and are not reachable through the Range/Worksheet/
#      Workbook COM surface, so they are left as-is here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Relabel the "CMS(e) " column (L2:L12) to "CMS(e)(2011)".
$ws.Range("L2:L12").Value = "CMS(e)(2011)"

# 2) Move the active selection to N7 (was N22), matching the saved
#    <selection activeCell="N7" sqref="N7"/> in the sheet view.
$ws.Range("N7").Select()
